$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain literal text
# (matching the original inlineStr formatting, e.g. trailing zeros like "1.000").
# Pre-format them as Text so Excel does not coerce the string into a Number.
$textForceCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D11",
    "D12",
    "D13",
    "D15",
    "D16",
    "D19",
    "D20",
    "D22",
    "D23",
    "D24",
    "D25",
    "D27",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D51"
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "26.522.33"
$ws.Range("E2").Value = "  +3.76%  "
$ws.Range("D3").Value = "1.740.18"
$ws.Range("E3").Value = "  +4.29%  "
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "245.37"
$ws.Range("E5").Value = "  +4.39%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.4809"
$ws.Range("E7").Value = "  +3.55%  "
$ws.Range("D8").Value = "0.2678"
$ws.Range("E8").Value = "  +4.04%  "
$ws.Range("D9").Value = "0.06243"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").Value = "1.737.50"
$ws.Range("E10").Value = "  +4.14%  "
$ws.Range("D11").Value = "0.07119"
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("D12").Value = "15.86"
$ws.Range("E12").Value = "  +8.30%  "
$ws.Range("D13").Value = "0.6208"
$ws.Range("E13").Value = "  +8.60%  "
$ws.Range("E14").Value = "  +4.65%  "
$ws.Range("D15").Value = "77.22"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D16").Value = "1.0000"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "26.530.80"
$ws.Range("E17").Value = "  +3.80%  "
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "0.000006890"
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("D20").Value = "11.77"
$ws.Range("E20").Value = "  +3.64%  "
$ws.Range("D21").Value = "1.962.12"
$ws.Range("E21").Value = "  +4.43%  "
$ws.Range("D22").Value = "4.591"
$ws.Range("E22").Value = "  +4.08%  "
$ws.Range("D23").Value = "8.912"
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("D24").Value = "5.354"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("D25").Value = "135.88"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  +3.61%  "
$ws.Range("D27").Value = "1.814"
$ws.Range("E27").Value = "  +6.00%  "
$ws.Range("E28").Value = "  +4.25%  "
$ws.Range("D29").Value = "106.88"
$ws.Range("E29").Value = "  +2.91%  "
$ws.Range("D30").Value = "4.002"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").Value = "3.743"
$ws.Range("E31").Value = "  +4.14%  "
$ws.Range("D32").Value = "0.07891"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("D33").Value = "0.04591"
$ws.Range("E33").Value = "  +6.40%  "
$ws.Range("D34").Value = "2.614"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.002"
$ws.Range("E35").Value = "  +6.27%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6381"
$ws.Range("E36").Value = "  +6.39%  "
$ws.Range("D37").Value = "0.9325"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("D38").Value = "112.17"
$ws.Range("E38").Value = "  +6.73%  "
$ws.Range("E39").Value = "  +8.27%  "
$ws.Range("D40").Value = "2.432"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.753"
$ws.Range("E42").Value = "  +14.14%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.01515"
$ws.Range("E43").Value = "  +3.56%  "
$ws.Range("D44").Value = "0.3916"
$ws.Range("E44").Value = "  +5.69%  "
$ws.Range("D45").Value = "6.971"
$ws.Range("E45").Value = "  +13.92%  "
$ws.Range("D46").Value = "0.1203"
$ws.Range("E46").Value = "  +8.39%  "
$ws.Range("D47").Value = "0.05332"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("D48").Value = "7.887"
$ws.Range("E48").Value = "  +4.00%  "
$ws.Range("D49").Value = "30.85"
$ws.Range("E49").Value = "  +3.76%  "
$ws.Range("E50").Value = "  +5.88%  "
$ws.Range("D51").Value = "0.3444"
$ws.Range("E51").Value = "  +4.41%  "
